$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 429, shifting existing rows 429-455 down to 430-456.
$ws.Rows("429:429").Insert()

# Populate the newly inserted row 429 with the new record's data.
$ws.Range("A429").Value = 6
$ws.Range("B429").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C429").Value = "Metropolitana"
$ws.Range("D429").Value = 44746
$ws.Range("D429").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E429").Value = 13
$ws.Range("F429").Value = 100112043
$ws.Range("G429").Value = "Pepino ensalada"
$ws.Range("H429").Value = "Sin especificar"
$ws.Range("I429").Value = "Primera"
$ws.Range("J429").Value = 550
$ws.Range("K429").Value = 15000
$ws.Range("L429").Value = 16000
$ws.Range("M429").Value = 15545
$ws.Range("N429").Value = '$/caja 50 unidades'
$ws.Range("O429").Value = "Región de Arica y Parinacota"
$ws.Range("P429").Value = 311
$ws.Range("Q429").Value = 50
$ws.Range("R429").Value = "Hortaliza"
